$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 395.58823
$ws.Range("J17").Value = 395.58823
$ws.Range("L17").Value = 1186.76469
$ws.Range("N17").Value = -1522.76469
$ws.Range("H129").Value = 47621036
$ws.Range("I129").Value = 2248.5
$ws.Range("J129").Value = 66668550
$ws.Range("K129").Value = 6745.5
$ws.Range("L129").Value = 200005650
$ws.Range("M129").Value = -1745.5
$ws.Range("N129").Value = -200015650
$ws.Range("H135").Value = 1874.0435
$ws.Range("I135").Value = 1036.375
$ws.Range("J135").Value = 3788.7144
$ws.Range("K135").Value = 9327.375
$ws.Range("L135").Value = 34098.4296
$ws.Range("M135").Value = -6792.375
$ws.Range("N135").Value = -39168.4296
$ws.Range("H137").Value = 4567.154
$ws.Range("I137").Value = 4331.467
$ws.Range("J137").Value = 4714.4585
$ws.Range("K137").Value = 12994.401
$ws.Range("L137").Value = 14143.3755
$ws.Range("M137").Value = -10444.401
$ws.Range("N137").Value = -19243.3755
$ws.Range("H138").Value = 9422.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 9422.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 28268.25
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -38548.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5145.25
$ws.Range("I2").Value = 3417
$ws.Range("J2").Value = 18971.25
$ws.Range("K2").Value = 3417
$ws.Range("L2").Value = 18971.25
$ws.Range("M2").Value = -3304
$ws.Range("N2").Value = -19197.25
$ws.Range("H32").Value = 1982.098
$ws.Range("I32").Value = 627.2619
$ws.Range("J32").Value = 8304.666999999999
$ws.Range("K32").Value = 627.2619
$ws.Range("L32").Value = 8304.666999999999
$ws.Range("M32").Value = -340.2619
$ws.Range("N32").Value = -8878.666999999999
$ws.Range("H45").Value = 1570.9546
$ws.Range("I45").Value = 1584.1904
$ws.Range("J45").Value = 1293
$ws.Range("K45").Value = 1584.1904
$ws.Range("L45").Value = 1293
$ws.Range("M45").Value = -1207.1904
$ws.Range("N45").Value = -2047
$ws.Range("H61").Value = 5332.788
$ws.Range("I61").Value = 5360.1787
$ws.Range("J61").Value = 5179.4
$ws.Range("K61").Value = 5360.1787
$ws.Range("L61").Value = 5179.4
$ws.Range("M61").Value = -5148.1787
$ws.Range("N61").Value = -5603.4
$ws.Range("H74").Value = 5499.3335
$ws.Range("I74").Value = 4500
$ws.Range("J74").Value = 5999
$ws.Range("K74").Value = 4500
$ws.Range("L74").Value = 5999
$ws.Range("M74").Value = -3626
$ws.Range("N74").Value = -7747
$ws.Range("H77").Value = 5499.3335
$ws.Range("I77").Value = 4500
$ws.Range("J77").Value = 5999
$ws.Range("K77").Value = 22500
$ws.Range("L77").Value = 29995
$ws.Range("M77").Value = -18132
$ws.Range("N77").Value = -38731
$ws.Range("H116").Value = 5145.25
$ws.Range("I116").Value = 3417
$ws.Range("J116").Value = 18971.25
$ws.Range("K116").Value = 3417
$ws.Range("L116").Value = 18971.25
$ws.Range("M116").Value = -1123
$ws.Range("N116").Value = -23559.25
$ws.Range("H132").Value = 70818.28999999999
$ws.Range("I132").Value = 1992.5714
$ws.Range("J132").Value = 139644
$ws.Range("K132").Value = 5977.7142
$ws.Range("L132").Value = 418932
$ws.Range("M132").Value = -3447.7142
$ws.Range("N132").Value = -423992
$ws.Range("H136").Value = 5332.788
$ws.Range("I136").Value = 5360.1787
$ws.Range("J136").Value = 5179.4
$ws.Range("K136").Value = 16080.5361
$ws.Range("L136").Value = 15538.2
$ws.Range("M136").Value = -13530.5361
$ws.Range("N136").Value = -20638.2

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5145.25
$ws.Range("I3").Value = 3417
$ws.Range("J3").Value = 18971.25
$ws.Range("K3").Value = 3417
$ws.Range("L3").Value = 18971.25
$ws.Range("M3").Value = -3303
$ws.Range("N3").Value = -19199.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4239.951
$ws.Range("I31").Value = 1076.5264
$ws.Range("J31").Value = 5194
$ws.Range("K31").Value = 1076.5264
$ws.Range("L31").Value = 5194
$ws.Range("M31").Value = -781.5264
$ws.Range("N31").Value = -5784
$ws.Range("H34").Value = 4239.951
$ws.Range("I34").Value = 1076.5264
$ws.Range("J34").Value = 5194
$ws.Range("K34").Value = 1076.5264
$ws.Range("L34").Value = 5194
$ws.Range("M34").Value = -874.5264
$ws.Range("N34").Value = -5598

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 160207.5
$ws.Range("I68").Value = 2911.75
$ws.Range("J68").Value = 212639.42
$ws.Range("K68").Value = 8735.25
$ws.Range("L68").Value = 637918.26
$ws.Range("M68").Value = -7924.25
$ws.Range("N68").Value = -639540.26
$ws.Range("H71").Value = 160207.5
$ws.Range("I71").Value = 2911.75
$ws.Range("J71").Value = 212639.42
$ws.Range("K71").Value = 26205.75
$ws.Range("L71").Value = 1913754.78
$ws.Range("M71").Value = -22149.75
$ws.Range("N71").Value = -1921866.78
$ws.Range("H127").Value = 4532.3335
$ws.Range("J127").Value = 4532.3335
$ws.Range("L127").Value = 13597.0005
$ws.Range("N127").Value = -23517.0005
$ws.Range("H137").Value = 9826
$ws.Range("I137").Value = 11591.2
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 34773.60000000001
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -29673.60000000001
$ws.Range("N137").Value = -13200

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2953.932
$ws.Range("I132").Value = 2869.4
$ws.Range("J132").Value = 3135.0715
$ws.Range("K132").Value = 8608.200000000001
$ws.Range("L132").Value = 9405.2145
$ws.Range("M132").Value = -6078.200000000001
$ws.Range("N132").Value = -14465.2145

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2370.366
$ws.Range("I132").Value = 2340.484
$ws.Range("J132").Value = 2463
$ws.Range("K132").Value = 7021.451999999999
$ws.Range("L132").Value = 7389
$ws.Range("M132").Value = -4491.451999999999
$ws.Range("N132").Value = -12449
$ws.Range("H136").Value = 32260948
$ws.Range("I136").Value = 1944.85
$ws.Range("J136").Value = 90913680
$ws.Range("K136").Value = 5834.549999999999
$ws.Range("L136").Value = 272741040
$ws.Range("M136").Value = -3284.549999999999
$ws.Range("N136").Value = -272746140

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2285.1277
$ws.Range("I136").Value = 1616.4857
$ws.Range("J136").Value = 4235.3335
$ws.Range("K136").Value = 4849.4571
$ws.Range("L136").Value = 12706.0005
$ws.Range("M136").Value = -2299.4571
$ws.Range("N136").Value = -17806.0005
